# "working on patient editor" -- add a new (blank) Sheet2 after Sheet1 and make
# it the active sheet, and tidy up a stray cell on Sheet1's row 4 (the label
# that was sitting in A4 belongs one column over, in B4; the red bottom-border
# formatting that used to live on B4 moves with it to C4).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: shift row 4's label from A4 -> B4, and its border format B4 -> C4 ---
$ws1.Range("B4").Copy()
$ws1.Range("C4").PasteSpecial(-4122)     # xlPasteFormats

$label = $ws1.Range("A4").Value2
$ws1.Range("A4").ClearContents()
$ws1.Range("B4").Value2 = $label
$ws1.Range("B4").Borders.LineStyle = 0   # B4 no longer carries the red border

# --- Sheet1: move the active selection to A5 ---
[void]$ws1.Range("A5").Select()

# --- Add Sheet2 right after Sheet1; it becomes the active/selected sheet ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
